$wb = $excel.ActiveWorkbook

# --- Processes sheet: "type" column moves from D to B, shifting old fuel/product
#     (B/C) one column to the right (new B=old D, new C=old B, new D=old C) ---
$wsProc = $wb.Worksheets.Item("Processes")
for ($r = 1; $r -le 9; $r++) {
    $oldB = $wsProc.Cells.Item($r, 2).Value2
    $oldC = $wsProc.Cells.Item($r, 3).Value2
    $oldD = $wsProc.Cells.Item($r, 4).Value2
    $wsProc.Cells.Item($r, 2).Value = $oldD
    $wsProc.Cells.Item($r, 3).Value = $oldB
    $wsProc.Cells.Item($r, 4).Value = $oldC
}
$wsProc.Range("B1:B9").Select() | Out-Null

# --- Exergy sheet: update stored selection ---
$wsExergy = $wb.Worksheets.Item("Exergy")
$wsExergy.Range("A2").Select() | Out-Null

# --- Flows sheet: A2:A16 become formulas pulling the key column from Exergy,
#     and the sheet becomes the active tab with A2:A16 selected ---
$wsFlows = $wb.Worksheets.Item("Flows")
$wsFlows.Activate() | Out-Null
for ($r = 2; $r -le 16; $r++) {
    $wsFlows.Cells.Item($r, 1).Formula = "=Exergy!A$r"
}
$wsFlows.Range("A2:A16").Select() | Out-Null
